$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 382 (2021-03-13): revised new-positive-cases count
$ws.Range("C382").Value = 43

# Row 383 (2021-03-14): revised new-positive-cases count
$ws.Range("C383").Value = 23

# Row 384 (2021-03-15): revised new-positive-cases count
$ws.Range("C384").Value = 66

# Row 385 (date 44271) was entirely blank (formulas returned "") - fill in the
# day's figures now that they are available.
$ws.Range("C385").Value = 11
$ws.Range("E385").Value = 6
$ws.Range("F385").Value = 5
$ws.Range("G385").Value = 31

# L385/M385 are formatted as Text (@). Assigning a numeric .Value directly to a
# Text-formatted cell makes Excel store it as a text string, which does not
# match how the rest of the column was populated (plain numeric cells with no
# "t" attribute). Temporarily drop to the default style to write a true
# number, then restore the original number format/borders by copying them
# from the cell directly above (which already carries the desired style) so
# no new style entries are introduced.
$ws.Range("L385").Style = "Normal"
$ws.Range("M385").Style = "Normal"
$ws.Range("L385").Value = 0
$ws.Range("M385").Value = 0
$ws.Range("L384").Copy()
$ws.Range("L385").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("M384").Copy()
$ws.Range("M385").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
